# Add a new "2022-Q3" sheet (fund holdings for the quarter) right after the
# "总计" (summary) sheet, and add the corresponding summary row on "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    below the header row, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 2.1

# Renumber the index column (A) for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计" with the fund
#    holdings detail for the quarter.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text in the source data (fund
# codes / percentages kept as strings), so force a text format before
# assigning the values to avoid Excel auto-converting them to numbers.
$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "011174"
$q3.Range("C2").Value = "中庚价值品质一年持有期混合"
$q3.Range("D2").Value = "66.33"
$q3.Range("E2").Value = "92.24"
$q3.Range("F2").Value = "3.06"
$q3.Range("G2").Value = "2.0297"
$q3.Range("H2").Value = 9

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "257050"
$q3.Range("C3").Value = "国联安主题驱动混合"
$q3.Range("D3").Value = "1.48"
$q3.Range("E3").Value = "93.02"
$q3.Range("F3").Value = "4.84"
$q3.Range("G3").Value = "0.0716"
$q3.Range("H3").Value = 4
